$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Poroto verde" series.
# It belongs chronologically at row 175 (pushing the existing rows
# 175-214 down to 176-215), so insert a blank row there first - Excel
# copies the formatting (incl. the date number format on column D)
# from the row above automatically.
$ws.Rows.Item(175).Insert()

$ws.Cells.Item(175, 1).Value = 5
$ws.Cells.Item(175, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(175, 3).Value = 'Maule'
$ws.Cells.Item(175, 4).Value2 = 44951
$ws.Cells.Item(175, 5).Value = 7
$ws.Cells.Item(175, 6).Value = 100112031
$ws.Cells.Item(175, 7).Value = 'Poroto verde'
$ws.Cells.Item(175, 8).Value = 'Sin especificar'
$ws.Cells.Item(175, 9).Value = 'Primera'
$ws.Cells.Item(175, 10).Value = 100
$ws.Cells.Item(175, 11).Value = 20000
$ws.Cells.Item(175, 12).Value = 20000
$ws.Cells.Item(175, 13).Value = 20000
$ws.Cells.Item(175, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(175, 15).Value = 'Región del Maule'
$ws.Cells.Item(175, 16).Value = 800
$ws.Cells.Item(175, 17).Value = 25
$ws.Cells.Item(175, 18).Value = 'Hortaliza'
